$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a text value into a cell (avoids Excel auto-converting
# numeric-looking strings like "381.22" into real numbers), then restore
# the cell to its original unstyled "General" appearance.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

$sub3 = [char]0x2083

Set-TextValue $ws.Range("D2") "51.060.26"
$ws.Range("E2").Value = "  -0.11%  "

Set-TextValue $ws.Range("D3") "2.957.67"
$ws.Range("E3").Value = "  +0.60%  "

$ws.Range("E4").Value = "  -0.02%  "

Set-TextValue $ws.Range("D5") "381.22"
$ws.Range("E5").Value = "  +1.31%  "

Set-TextValue $ws.Range("D6") "102.22"
$ws.Range("E6").Value = "  -0.16%  "

Set-TextValue $ws.Range("D7") "0.545"
$ws.Range("E7").Value = "  +1.88%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("E9").Value = "  +0.75%  "

Set-TextValue $ws.Range("D10") "36.51"
$ws.Range("E10").Value = "  -0.79%  "

$ws.Range("E11").Value = "  -0.73%  "

$ws.Range("E12").Value = "  +1.73%  "

Set-TextValue $ws.Range("D13") "12.48"
$ws.Range("E13").Value = "  +75.29%  "

$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Range("D14") "18.39"
$ws.Range("E14").Value = "  +2.55%  "

$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue $ws.Range("D15") "3.415.31"
$ws.Range("E15").Value = "  +0.36%  "

Set-TextValue $ws.Range("D16") "7.73"
$ws.Range("E16").Value = "  +5.24%  "

Set-TextValue $ws.Range("D17") "2.964.76"
$ws.Range("E17").Value = "  +1.05%  "

$ws.Range("E18").Value = "  +3.51%  "

Set-TextValue $ws.Range("D19") "51.108.03"
$ws.Range("E19").Value = "  +0.08%  "

$ws.Range("E20").Value = "  -2.31%  "

Set-TextValue $ws.Range("D21") "12.39"
$ws.Range("E21").Value = "  -1.40%  "

$d22val = "0.0{0}0959" -f $sub3
Set-TextValue $ws.Range("D22") $d22val
$ws.Range("E22").Value = "  +0.35%  "

Set-TextValue $ws.Range("D23") "3.33"
$ws.Range("E23").Value = "  +16.28%  "

Set-TextValue $ws.Range("D24") "270.10"
$ws.Range("E24").Value = "  +2.58%  "

Set-TextValue $ws.Range("D25") "69.72"
$ws.Range("E25").Value = "  +2.19%  "

Set-TextValue $ws.Range("D26") "7.92"
$ws.Range("E26").Value = "  -3.33%  "

Set-TextValue $ws.Range("D27") "0.998"
$ws.Range("E27").Value = "  -0.16%  "

$ws.Range("E28").Value = "  -0.23%  "

Set-TextValue $ws.Range("D29") "25.90"
$ws.Range("E29").Value = "  +0.89%  "

Set-TextValue $ws.Range("D30") "7.02"
$ws.Range("E30").Value = "  -10.88%  "

$ws.Range("E31").Value = "  -3.33%  "

Set-TextValue $ws.Range("D32") "10.43"
$ws.Range("E32").Value = "  +5.82%  "

Set-TextValue $ws.Range("D33") "51.18"
$ws.Range("E33").Value = "  +0.98%  "

Set-TextValue $ws.Range("D34") "34.26"

$ws.Range("E35").Value = "  +2.04%  "

$ws.Range("E36").Value = "  -3.92%  "

$ws.Range("E37").Value = "  -0.05%  "

$ws.Range("E38").Value = "  +9.21%  "

$ws.Range("E39").Value = "  +2.16%  "

$ws.Range("E40").Value = "  +1.80%  "

$ws.Range("E41").Value = "  +2.79%  "

Set-TextValue $ws.Range("D42") "2.50"
$ws.Range("E42").Value = "  -2.83%  "

Set-TextValue $ws.Range("D43") "124.61"
$ws.Range("E43").Value = "  +2.45%  "

Set-TextValue $ws.Range("D44") "21.73"
$ws.Range("E44").Value = "  +3.00%  "

Set-TextValue $ws.Range("D45") "3.56"
$ws.Range("E45").Value = "  +10.49%  "

Set-TextValue $ws.Range("D46") "2.073.65"
$ws.Range("E46").Value = "  +3.59%  "

Set-TextValue $ws.Range("D47") "2.02"
$ws.Range("E47").Value = "  -1.33%  "

$ws.Range("E48").Value = "  +0.10%  "

Set-TextValue $ws.Range("D49") "0.261"
$ws.Range("E49").Value = "  -4.89%  "

$ws.Range("E50").Value = "  -7.28%  "

$ws.Range("E51").Value = "  +6.58%  "

